$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.32
$ws.Range("F2").Value = 0.73

$ws.Range("D3").Value = 1.43
$ws.Range("E3").Value = 1.26

$ws.Range("C4").Value = 1.4
$ws.Range("F4").Value = 1.09

$ws.Range("C5").Value = 1.38
$ws.Range("E5").Value = 1.21
$ws.Range("F5").Value = 1.05

$ws.Range("D6").Value = 1.52
$ws.Range("E6").Value = 1.31
$ws.Range("F6").Value = 1.15
